# The document contains a single inline picture ("Picture 9") sitting in
# its own paragraph, right after the closing paragraph mark of the
# preceding body text. The edit removes that picture (its containing
# <w:r> run) entirely, leaving the now-empty paragraph (just its
# paragraph properties) followed by the blank paragraph that already
# existed after it.

$d = $word.ActiveDocument

# Walk the inline shapes back-to-front and drop any picture shapes —
# robust even if more than one were ever present.
for ($i = $d.InlineShapes.Count; $i -ge 1; $i--) {
    $shape = $d.InlineShapes.Item($i)
    $shape.Delete()
}
